# khl_stats_1369_ext.xlsx update (2025-11-19 refresh)
# - Matches_SOG: append two new completed matches (rows 494-495)
# - Shots_HA / Shots_Summary: bump as_of_utc to the new snapshot and refresh
#   the four teams (Драконы, Локомотив, Салават Юлаев, ХК Сочи) whose
#   cumulative shot counts moved because of those matches
# - Meta_ext: bump as_of_utc + build_version

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Matches_SOG: append the two newly completed games
# ---------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$wsMatches.Cells.Item(494, 1).Value = "'897791"
$wsMatches.Cells.Item(494, 2).Value = "2025-11-19T19:00:00"
$wsMatches.Cells.Item(494, 3).Value = "Локомотив"
$wsMatches.Cells.Item(494, 4).Value = "Драконы"
$wsMatches.Cells.Item(494, 5).Value = 49
$wsMatches.Cells.Item(494, 6).Value = 14
$wsMatches.Cells.Item(494, 7).Value = "khl_text"

$wsMatches.Cells.Item(495, 1).Value = "'897792"
$wsMatches.Cells.Item(495, 2).Value = "2025-11-19T19:30:00"
$wsMatches.Cells.Item(495, 3).Value = "ХК Сочи"
$wsMatches.Cells.Item(495, 4).Value = "Салават Юлаев"
$wsMatches.Cells.Item(495, 5).Value = 31
$wsMatches.Cells.Item(495, 6).Value = 30
$wsMatches.Cells.Item(495, 7).Value = "khl_text"

# ---------------------------------------------------------------------
# Shots_HA: refresh as_of_utc for every team, then update the four teams
# whose totals changed as a result of the new games
# ---------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

$wsHA.Range("D2:D23").Value = "2025-11-19T19:30:00Z"

# Row 10 - Драконы (away in game 897791)
$wsHA.Cells.Item(10, 6).Value = 24
$wsHA.Cells.Item(10, 11).Value = 649
$wsHA.Cells.Item(10, 12).Value = 869
$wsHA.Cells.Item(10, 13).Value = 27
$wsHA.Cells.Item(10, 14).Value = 36.2

# Row 12 - Локомотив (home in game 897791)
$wsHA.Cells.Item(12, 5).Value = 21
$wsHA.Cells.Item(12, 7).Value = 667
$wsHA.Cells.Item(12, 8).Value = 561
$wsHA.Cells.Item(12, 9).Value = 31.8
$wsHA.Cells.Item(12, 10).Value = 26.7

# Row 16 - Салават Юлаев (away in game 897792)
$wsHA.Cells.Item(16, 6).Value = 28
$wsHA.Cells.Item(16, 11).Value = 781
$wsHA.Cells.Item(16, 12).Value = 825
$wsHA.Cells.Item(16, 13).Value = 27.9
$wsHA.Cells.Item(16, 14).Value = 29.5

# Row 22 - ХК Сочи (home in game 897792)
$wsHA.Cells.Item(22, 5).Value = 22
$wsHA.Cells.Item(22, 7).Value = 643
$wsHA.Cells.Item(22, 8).Value = 718
$wsHA.Cells.Item(22, 9).Value = 29.2
$wsHA.Cells.Item(22, 10).Value = 32.6

# ---------------------------------------------------------------------
# Shots_Summary: same refresh pattern
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

$wsSummary.Range("D2:D23").Value = "2025-11-19T19:30:00Z"

# Row 10 - Драконы
$wsSummary.Cells.Item(10, 5).Value = 44
$wsSummary.Cells.Item(10, 6).Value = 1220
$wsSummary.Cells.Item(10, 7).Value = 1567
$wsSummary.Cells.Item(10, 8).Value = 27.7
$wsSummary.Cells.Item(10, 9).Value = 35.6

# Row 12 - Локомотив
$wsSummary.Cells.Item(12, 5).Value = 48
$wsSummary.Cells.Item(12, 6).Value = 1506
$wsSummary.Cells.Item(12, 7).Value = 1227
$wsSummary.Cells.Item(12, 8).Value = 31.4
$wsSummary.Cells.Item(12, 9).Value = 25.6

# Row 16 - Салават Юлаев
$wsSummary.Cells.Item(16, 5).Value = 45
$wsSummary.Cells.Item(16, 6).Value = 1246
$wsSummary.Cells.Item(16, 7).Value = 1294
$wsSummary.Cells.Item(16, 8).Value = 27.7
$wsSummary.Cells.Item(16, 9).Value = 28.8

# Row 22 - ХК Сочи
$wsSummary.Cells.Item(22, 5).Value = 42
$wsSummary.Cells.Item(22, 6).Value = 1159
$wsSummary.Cells.Item(22, 7).Value = 1463
$wsSummary.Cells.Item(22, 8).Value = 27.6
$wsSummary.Cells.Item(22, 9).Value = 34.8

# ---------------------------------------------------------------------
# Meta_ext: bump snapshot timestamp + build number
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-11-19T19:30:00Z"
$wsMeta.Cells.Item(2, 4).Value = 77

Write-Output "khl_stats_1369_ext.xlsx updated for 2025-11-19 snapshot"
